$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.775841999999999
$ws.Range("H2").Value = 8.327525999999999
$ws.Range("I2").Value = 0.0624750527258915
$ws.Range("J2").Value = 0.0624750527258915
$ws.Range("M2").Value = 0.8366046666666667
$ws.Range("N2").Value = 2.509814
$ws.Range("O2").Value = 0.08025679986157715
$ws.Range("P2").Value = 0.08025679986157715
$ws.Range("Q2").Value = 2.322282371129333
$ws.Range("R2").Value = 20.900541340164
$ws.Range("S2").Value = 0.005014047802963354
$ws.Range("T2").Value = 0.005014047802963354
$ws.Range("G3").Value = 2.775841999999999
$ws.Range("H3").Value = 8.327525999999999
$ws.Range("I3").Value = 0.0624750527258915
$ws.Range("J3").Value = 0.0624750527258915
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("O3").Value = 0.7616247559221037
$ws.Range("P3").Value = 0.7616247559221038
$ws.Range("Q3").Value = 22.03810452378066
$ws.Range("R3").Value = 198.342940714026
$ws.Range("S3").Value = 0.04758254678357766
$ws.Range("T3").Value = 0.04758254678357768
$ws.Range("G4").Value = 2.775841999999999
$ws.Range("H4").Value = 8.327525999999999
$ws.Range("I4").Value = 0.0624750527258915
$ws.Range("J4").Value = 0.0624750527258915
$ws.Range("M4").Value = 1.648242
$ws.Range("N4").Value = 4.944726
$ws.Range("O4").Value = 0.1581184442163192
$ws.Range("P4").Value = 0.1581184442163192
$ws.Range("Q4").Value = 4.575259369763999
$ws.Range("R4").Value = 41.177334327876
$ws.Range("S4").Value = 0.009878458139350474
$ws.Range("T4").Value = 0.009878458139350477
$ws.Range("I5").Value = 0.2652892219050753
$ws.Range("J5").Value = 0.2652892219050753
$ws.Range("M5").Value = 0.8366046666666667
$ws.Range("N5").Value = 2.509814
$ws.Range("O5").Value = 0.08025679986157715
$ws.Range("P5").Value = 0.08025679986157715
$ws.Range("Q5").Value = 9.861159877428223
$ws.Range("R5").Value = 88.75043889685401
$ws.Range("S5").Value = 0.02129126398786915
$ws.Range("T5").Value = 0.02129126398786915
$ws.Range("I6").Value = 0.2652892219050753
$ws.Range("J6").Value = 0.2652892219050753
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("O6").Value = 0.7616247559221037
$ws.Range("P6").Value = 0.7616247559221038
$ws.Range("Q6").Value = 93.58089903545678
$ws.Range("R6").Value = 842.2280913191109
$ws.Range("S6").Value = 0.2020508388822177
$ws.Range("T6").Value = 0.2020508388822178
$ws.Range("I7").Value = 0.2652892219050753
$ws.Range("J7").Value = 0.2652892219050753
$ws.Range("M7").Value = 1.648242
$ws.Range("N7").Value = 4.944726
$ws.Range("O7").Value = 0.1581184442163192
$ws.Range("P7").Value = 0.1581184442163192
$ws.Range("Q7").Value = 19.428026792454
$ws.Range("R7").Value = 174.852241132086
$ws.Range("S7").Value = 0.04194711903498836
$ws.Range("T7").Value = 0.04194711903498837
$ws.Range("G8").Value = 29.86824466666667
$ws.Range("H8").Value = 89.60473400000001
$ws.Range("I8").Value = 0.6722357253690333
$ws.Range("J8").Value = 0.6722357253690333
$ws.Range("M8").Value = 0.8366046666666667
$ws.Range("N8").Value = 2.509814
$ws.Range("O8").Value = 0.08025679986157715
$ws.Range("P8").Value = 0.08025679986157715
$ws.Range("Q8").Value = 24.98791287327511
$ws.Range("R8").Value = 224.891215859476
$ws.Range("S8").Value = 0.05395148807074465
$ws.Range("T8").Value = 0.05395148807074465
$ws.Range("G9").Value = 29.86824466666667
$ws.Range("H9").Value = 89.60473400000001
$ws.Range("I9").Value = 0.6722357253690333
$ws.Range("J9").Value = 0.6722357253690333
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("O9").Value = 0.7616247559221037
$ws.Range("P9").Value = 0.7616247559221038
$ws.Range("Q9").Value = 237.1314714259149
$ws.Range("R9").Value = 2134.183242833234
$ws.Range("S9").Value = 0.5119913702563083
$ws.Range("T9").Value = 0.5119913702563084
$ws.Range("G10").Value = 29.86824466666667
$ws.Range("H10").Value = 89.60473400000001
$ws.Range("I10").Value = 0.6722357253690333
$ws.Range("J10").Value = 0.6722357253690333
$ws.Range("M10").Value = 1.648242
$ws.Range("N10").Value = 4.944726
$ws.Range("O10").Value = 0.1581184442163192
$ws.Range("P10").Value = 0.1581184442163192
$ws.Range("Q10").Value = 49.230095325876
$ws.Range("R10").Value = 443.0708579328841
$ws.Range("S10").Value = 0.1062928670419803
$ws.Range("T10").Value = 0.1062928670419804
